$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$m.Theme.Name = "Office Theme"
